# Update "想去人数" (F column) counts across sheets as per the data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 801
$ws1.Range("F5").Value = 1062
$ws1.Range("F8").Value = 208
$ws1.Range("F9").Value = 383
$ws1.Range("F10").Value = 3
$ws1.Range("F12").Value = 498
$ws1.Range("F15").Value = 12471
$ws1.Range("F16").Value = 135
$ws1.Range("F17").Value = 5495

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 122

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 801
$ws4.Range("F4").Value = 122
$ws4.Range("F7").Value = 1062
$ws4.Range("F10").Value = 208
$ws4.Range("F11").Value = 383
$ws4.Range("F12").Value = 3
$ws4.Range("F14").Value = 498
$ws4.Range("F17").Value = 12471
$ws4.Range("F19").Value = 135
$ws4.Range("F20").Value = 5495
